$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bulk cell value updates (odds/match data refresh) ---
$ws.Range("AB156").Value = 0.8999999999999999
$ws.Range("B156").Value = 7211641
$ws.Range("E156").Value = 'Sport Huancayo'
$ws.Range("F156").Value = 'Deportivo Municipal'
$ws.Range("G156").Value = 2
$ws.Range("H156").Value = 0
$ws.Range("I156").Value = 'H'
$ws.Range("J156").Value = 1.125
$ws.Range("K156").Value = 7
$ws.Range("L156").Value = 17
$ws.Range("M156").Value = 1.166
$ws.Range("N156").Value = 6.5
$ws.Range("O156").Value = 12
$ws.Range("P156").Value = -2
$ws.Range("Q156").Value = 1.775
$ws.Range("R156").Value = 2.025
$ws.Range("S156").Value = 3.5
$ws.Range("T156").Value = 1.9
$ws.Range("U156").Value = 1.9
$ws.Range("V156").Value = 0.1659999999999999
$ws.Range("W156").Value = -1
$ws.Range("Y156").Value = 0
$ws.Range("Z156").Value = 0
$ws.Range("AB157").Value = 0.9750000000000001
$ws.Range("B157").Value = 7211640
$ws.Range("E157").Value = 'UTC Cajamarca'
$ws.Range("F157").Value = 'Sport Boys'
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 1
$ws.Range("I157").Value = 'D'
$ws.Range("J157").Value = 1.615
$ws.Range("K157").Value = 3.75
$ws.Range("L157").Value = 5
$ws.Range("M157").Value = 1.5
$ws.Range("N157").Value = 4.2
$ws.Range("O157").Value = 6.5
$ws.Range("P157").Value = -1
$ws.Range("Q157").Value = 1.8
$ws.Range("R157").Value = 2.05
$ws.Range("S157").Value = 2.5
$ws.Range("T157").Value = 1.875
$ws.Range("U157").Value = 1.975
$ws.Range("V157").Value = -1
$ws.Range("W157").Value = 3.2
$ws.Range("Y157").Value = -1
$ws.Range("Z157").Value = 1.05
$ws.Range("AB184").Value = 0.9750000000000001
$ws.Range("B184").Value = 7384629
$ws.Range("E184").Value = 'Deportivo Garcilaso'
$ws.Range("F184").Value = 'Alianza Lima'
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 1
$ws.Range("I184").Value = 'A'
$ws.Range("J184").Value = 2.625
$ws.Range("K184").Value = 3.3
$ws.Range("L184").Value = 2.5
$ws.Range("M184").Value = 2.7
$ws.Range("N184").Value = 3.4
$ws.Range("O184").Value = 2.375
$ws.Range("P184").Value = 0
$ws.Range("Q184").Value = 2.025
$ws.Range("R184").Value = 1.775
$ws.Range("S184").Value = 2.25
$ws.Range("T184").Value = 1.825
$ws.Range("U184").Value = 1.975
$ws.Range("V184").Value = -1
$ws.Range("X184").Value = 1.375
$ws.Range("Y184").Value = -1
$ws.Range("Z184").Value = 0.7749999999999999
$ws.Range("AA185").Value = -1
$ws.Range("AB185").Value = 0.95
$ws.Range("B185").Value = 7384627
$ws.Range("E185").Value = 'Universitario de Deportes'
$ws.Range("F185").Value = 'Sport Huancayo'
$ws.Range("G185").Value = 2
$ws.Range("J185").Value = 1.25
$ws.Range("L185").Value = 12
$ws.Range("M185").Value = 1.181
$ws.Range("N185").Value = 6
$ws.Range("P185").Value = -1.75
$ws.Range("Q185").Value = 1.8
$ws.Range("R185").Value = 2
$ws.Range("S185").Value = 2.75
$ws.Range("T185").Value = 1.85
$ws.Range("U185").Value = 1.95
$ws.Range("V185").Value = 0.181
$ws.Range("Y185").Value = 0.4
$ws.Range("Z185").Value = -0.5
$ws.Range("AA187").Value = 0.4875
$ws.Range("AB187").Value = -0.5
$ws.Range("B187").Value = 7384628
$ws.Range("E187").Value = 'Deportivo Binacional'
$ws.Range("F187").Value = 'FBC Melgar'
$ws.Range("G187").Value = 1
$ws.Range("H187").Value = 2
$ws.Range("J187").Value = 2.75
$ws.Range("L187").Value = 2.375
$ws.Range("M187").Value = 3.3
$ws.Range("N187").Value = 3.6
$ws.Range("O187").Value = 2
$ws.Range("P187").Value = 0.5
$ws.Range("Q187").Value = 1.8
$ws.Range("R187").Value = 2
$ws.Range("S187").Value = 2.75
$ws.Range("T187").Value = 1.975
$ws.Range("U187").Value = 1.875
$ws.Range("X187").Value = 1
$ws.Range("Z187").Value = 1
$ws.Range("AA188").Value = -0.5
$ws.Range("AB188").Value = 0.4
$ws.Range("B188").Value = 7384626
$ws.Range("E188").Value = 'Sporting Cristal'
$ws.Range("F188").Value = 'Alianza Atletico'
$ws.Range("G188").Value = 3
$ws.Range("H188").Value = 0
$ws.Range("I188").Value = 'H'
$ws.Range("J188").Value = 1.3
$ws.Range("K188").Value = 5
$ws.Range("L188").Value = 9
$ws.Range("M188").Value = 1.166
$ws.Range("N188").Value = 6.5
$ws.Range("O188").Value = 13
$ws.Range("P188").Value = -2
$ws.Range("Q188").Value = 1.85
$ws.Range("R188").Value = 1.95
$ws.Range("S188").Value = 3.25
$ws.Range("T188").Value = 2
$ws.Range("U188").Value = 1.8
$ws.Range("V188").Value = 0.1659999999999999
$ws.Range("X188").Value = -1
$ws.Range("Y188").Value = 0.8500000000000001
$ws.Range("Z188").Value = -1
$ws.Range("D336").Value = 45436.9375
$ws.Range("E336").Value = 'Cesar Vallejo'
$ws.Range("F336").Value = 'Atletico Grau'
$ws.Range("J336").Value = 1.95
$ws.Range("K336").Value = 3.2
$ws.Range("L336").Value = 4
$ws.Range("M336").Value = 1.909
$ws.Range("N336").Value = 3.2
$ws.Range("O336").Value = 4.1
$ws.Range("P336").Value = -0.5
$ws.Range("Q336").Value = 2
$ws.Range("R336").Value = 1.85
$ws.Range("S336").Value = 2.25
$ws.Range("T336").Value = 1.825
$ws.Range("U336").Value = 2.025
$ws.Range("D337").Value = 45437.54166666666
$ws.Range("E337").Value = 'Sport Huancayo'
$ws.Range("F337").Value = 'UTC Cajamarca'
$ws.Range("J337").Value = 1.533
$ws.Range("K337").Value = 4
$ws.Range("L337").Value = 6
$ws.Range("M337").Value = 1.65
$ws.Range("N337").Value = 4
$ws.Range("O337").Value = 4.75
$ws.Range("P337").Value = -0.75
$ws.Range("Q337").Value = 1.875
$ws.Range("R337").Value = 1.975
$ws.Range("S337").Value = 2.5
$ws.Range("T337").Value = 2
$ws.Range("U337").Value = 1.85
$ws.Range("T338").Value = 1.875
$ws.Range("U338").Value = 1.975
$ws.Range("E339").Value = 'Comerciantes Unidos'
$ws.Range("F339").Value = 'Sporting Cristal'
$ws.Range("J339").Value = 9
$ws.Range("K339").Value = 5
$ws.Range("L339").Value = 1.333
$ws.Range("M339").Value = 6.5
$ws.Range("N339").Value = 5
$ws.Range("O339").Value = 1.42
$ws.Range("P339").Value = 1.25
$ws.Range("Q339").Value = 1.975
$ws.Range("R339").Value = 1.875
$ws.Range("S339").Value = 3
$ws.Range("T339").Value = 1.875
$ws.Range("U339").Value = 1.975
$ws.Range("D340").Value = 45437.70833333334
$ws.Range("E340").Value = 'Deportivo Garcilaso'
$ws.Range("F340").Value = 'FBC Melgar'
$ws.Range("J340").Value = 2.625
$ws.Range("L340").Value = 2.75
$ws.Range("M340").Value = 3.9
$ws.Range("N340").Value = 3.25
$ws.Range("O340").Value = 1.95
$ws.Range("P340").Value = 0.5
$ws.Range("Q340").Value = 1.825
$ws.Range("R340").Value = 2.025
$ws.Range("S340").Value = 2.5
$ws.Range("T340").Value = 1.95
$ws.Range("U340").Value = 1.9
$ws.Range("M342").Value = 2.45
$ws.Range("O342").Value = 2.8
$ws.Range("P342").Value = 0
$ws.Range("Q342").Value = 1.8
$ws.Range("R342").Value = 2.05
$ws.Range("Q343").Value = 2.025
$ws.Range("R343").Value = 1.825
$ws.Range("T343").Value = 1.95
$ws.Range("U343").Value = 1.9

# --- "id" column for the most-recent batch of rows (335-343) is stored as
# text (shared-string) in the source data. Force text formatting on just the
# four rows whose id actually changes so the values stay text instead of
# being auto-coerced to numbers, matching the source refresh semantics.
foreach ($r in 336,337,339,340) {
    $ws.Cells.Item($r, 2).NumberFormat = "@"
}
$ws.Cells.Item(336, 2).Value = "8240874"
$ws.Cells.Item(337, 2).Value = "8240871"
$ws.Cells.Item(339, 2).Value = "8240872"
$ws.Cells.Item(340, 2).Value = "8240870"
